$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the values up: A1 <- "Print graph", A2 <- "Automatenmodus"
$ws.Range("A1").Value = "Print graph"
$ws.Range("A2").Value = "Automatenmodus"

# Remove the now-obsolete third row entirely
$ws.Rows.Item(3).Delete()

# Update the selection to match the new data range
$ws.Range("A1:A2").Select()
